$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.599.32"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.612.05"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'590.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'149.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").Value = "'0.110"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'5.77"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'27.68"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "3.073.95"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "63.406.85"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "'0.0000157"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "2.612.83"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "'12.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'4.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "'345.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'6.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'66.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'1.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").Value = "'9.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'1.66"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "'8.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").Value = "'547.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "0.0₃0867"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'5.34"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "'6.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'166.18"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'19.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'165.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "'4.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'23.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.76%  "
$ws.Range("D45").Value = "'0.0581"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").Value = "'2.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.50%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'0.0251"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'0.0960"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "'19.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "0.0₆0230"
$ws.Range("E51").Value = "  +15.67%  "
